$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new country code row (row iterator style: next empty row in column A)
$ws.Cells.Item(3, 1).Value = "CHN"

# Move the selection (as left by the editing session)
$ws.Range("C5").Select()
